# Update date and regenerate the multiplication-table answer values
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-12-14 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-15 Monday", 2) | Out-Null
$d.Content.Find.Execute("65×15=975", $true, $false, $false, $false, $false, $true, 1, $false, "15×33=495", 2) | Out-Null
$d.Content.Find.Execute("36×12=432", $true, $false, $false, $false, $false, $true, 1, $false, "57×14=798", 2) | Out-Null
$d.Content.Find.Execute("37×17=629", $true, $false, $false, $false, $false, $true, 1, $false, "99×14=1386", 2) | Out-Null
$d.Content.Find.Execute("15×93=1395", $true, $false, $false, $false, $false, $true, 1, $false, "14×74=1036", 2) | Out-Null
$d.Content.Find.Execute("51×81=4131", $true, $false, $false, $false, $false, $true, 1, $false, "12×62=744", 2) | Out-Null
$d.Content.Find.Execute("16×40=640", $true, $false, $false, $false, $false, $true, 1, $false, "18×81=1458", 2) | Out-Null
$d.Content.Find.Execute("19×33=627", $true, $false, $false, $false, $false, $true, 1, $false, "90×72=6480", 2) | Out-Null
$d.Content.Find.Execute("56×65=3640", $true, $false, $false, $false, $false, $true, 1, $false, "46×81=3726", 2) | Out-Null
$d.Content.Find.Execute("35×98=3430", $true, $false, $false, $false, $false, $true, 1, $false, "11×59=649", 2) | Out-Null
$d.Content.Find.Execute("71×44=3124", $true, $false, $false, $false, $false, $true, 1, $false, "67×26=1742", 2) | Out-Null
$d.Content.Find.Execute("65×63=4095", $true, $false, $false, $false, $false, $true, 1, $false, "66×77=5082", 2) | Out-Null
$d.Content.Find.Execute("93×19=1767", $true, $false, $false, $false, $false, $true, 1, $false, "19×42=798", 2) | Out-Null
$d.Content.Find.Execute("90×79=7110", $true, $false, $false, $false, $false, $true, 1, $false, "25×76=1900", 2) | Out-Null
$d.Content.Find.Execute("95×21=1995", $true, $false, $false, $false, $false, $true, 1, $false, "85×12=1020", 2) | Out-Null
$d.Content.Find.Execute("94×43=4042", $true, $false, $false, $false, $false, $true, 1, $false, "33×76=2508", 2) | Out-Null
$d.Content.Find.Execute("56×39=2184", $true, $false, $false, $false, $false, $true, 1, $false, "58×61=3538", 2) | Out-Null
$d.Content.Find.Execute("59×66=3894", $true, $false, $false, $false, $false, $true, 1, $false, "68×45=3060", 2) | Out-Null
$d.Content.Find.Execute("60×54=3240", $true, $false, $false, $false, $false, $true, 1, $false, "92×23=2116", 2) | Out-Null
$d.Content.Find.Execute("34×99=3366", $true, $false, $false, $false, $false, $true, 1, $false, "82×17=1394", 2) | Out-Null
$d.Content.Find.Execute("61×33=2013", $true, $false, $false, $false, $false, $true, 1, $false, "37×87=3219", 2) | Out-Null
$d.Content.Find.Execute("34×57=1938", $true, $false, $false, $false, $false, $true, 1, $false, "80×22=1760", 2) | Out-Null
$d.Content.Find.Execute("15×14=210", $true, $false, $false, $false, $false, $true, 1, $false, "50×91=4550", 2) | Out-Null
$d.Content.Find.Execute("13×26=338", $true, $false, $false, $false, $false, $true, 1, $false, "27×99=2673", 2) | Out-Null
$d.Content.Find.Execute("46×98=4508", $true, $false, $false, $false, $false, $true, 1, $false, "30×47=1410", 2) | Out-Null
$d.Content.Find.Execute("72×34=2448", $true, $false, $false, $false, $false, $true, 1, $false, "39×34=1326", 2) | Out-Null
